# Updates the crypto price/volume table (rows 2-51) with refreshed values,
# including a couple of row swaps (XRP/Dogecoin, Aptos/Litecoin, Filecoin/Aave)
# as captured by the upstream GitHub Actions data refresh.
#
# For cells in column D whose new value is a plain decimal number (e.g. "218.47"),
# we briefly force a Text number format before assigning so Excel stores the
# literal digit string instead of silently converting it to a numeric value
# (which would also mangle values like "0.0000254" into scientific notation,
# or "1.30" into "1.3"). The original cell style is restored immediately after,
# so no visible/persisted formatting change is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '91.295.98'
$ws.Range("E2").Value = '  +3.71%  '
$ws.Range("D3").Value = '3.114.89'
$ws.Range("E3").Value = '  +1.96%  '
$ws.Range("E4").Value = '  +0.00%  '
$c = $ws.Range("D5")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '218.47'
$c.Style = $origStyle
$ws.Range("E5").Value = '  +3.73%  '
$c = $ws.Range("D6")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '622.07'
$c.Style = $origStyle
$ws.Range("E6").Value = '  +0.66%  '
$ws.Range("B7").Value = 'XRP'
$ws.Range("C7").Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$c = $ws.Range("D7")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '0.986'
$c.Style = $origStyle
$ws.Range("E7").Value = '  +24.02%  '
$ws.Range("B8").Value = 'Dogecoin'
$ws.Range("C8").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$c = $ws.Range("D8")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '0.379'
$c.Style = $origStyle
$ws.Range("E8").Value = '  +2.28%  '
$ws.Range("E9").Value = '  +0.03%  '
$ws.Range("D10").Value = '3.110.27'
$ws.Range("E10").Value = '  +1.94%  '
$c = $ws.Range("D11")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '0.714'
$c.Style = $origStyle
$ws.Range("E11").Value = '  +20.45%  '
$ws.Range("E12").Value = '  +6.93%  '
$c = $ws.Range("D13")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '0.0000254'
$c.Style = $origStyle
$ws.Range("E13").Value = '  +6.96%  '
$c = $ws.Range("D14")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '34.57'
$c.Style = $origStyle
$ws.Range("E14").Value = '  +8.18%  '
$c = $ws.Range("D15")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '5.41'
$c.Style = $origStyle
$ws.Range("E15").Value = '  +2.86%  '
$ws.Range("D16").Value = '91.080.38'
$ws.Range("E16").Value = '  +3.75%  '
$ws.Range("D17").Value = '3.681.25'
$ws.Range("E17").Value = '  +1.94%  '
$ws.Range("D18").Value = '3.126.27'
$ws.Range("E18").Value = '  +2.46%  '
$c = $ws.Range("D19")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '3.77'
$c.Style = $origStyle
$ws.Range("E19").Value = '  +14.41%  '
$c = $ws.Range("D20")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '0.0000218'
$c.Style = $origStyle
$ws.Range("E20").Value = '  +6.46%  '
$c = $ws.Range("D21")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '14.06'
$c.Style = $origStyle
$ws.Range("E21").Value = '  +5.78%  '
$c = $ws.Range("D22")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '435.19'
$c.Style = $origStyle
$ws.Range("E22").Value = '  +3.83%  '
$c = $ws.Range("D23")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '8.83'
$c.Style = $origStyle
$ws.Range("E23").Value = '  +8.54%  '
$ws.Range("E24").Value = '  +5.72%  '
$c = $ws.Range("D25")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '6.15'
$c.Style = $origStyle
$ws.Range("E25").Value = '  +12.46%  '
$ws.Range("B26").Value = 'Litecoin'
$ws.Range("C26").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$c = $ws.Range("D26")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '86.45'
$c.Style = $origStyle
$ws.Range("E26").Value = '  +6.00%  '
$ws.Range("B27").Value = 'Aptos'
$ws.Range("C27").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$c = $ws.Range("D27")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '12.21'
$c.Style = $origStyle
$ws.Range("E27").Value = '  +3.77%  '
$ws.Range("D28").Value = '3.274.00'
$ws.Range("E28").Value = '  +2.13%  '
$ws.Range("E29").Value = '  +0.04%  '
$ws.Range("E30").Value = '  -4.65%  '
$c = $ws.Range("D31")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '9.08'
$c.Style = $origStyle
$ws.Range("E31").Value = '  +13.12%  '
$c = $ws.Range("D32")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '524.46'
$c.Style = $origStyle
$ws.Range("E32").Value = '  +3.68%  '
$ws.Range("E33").Value = '  -17.08%  '
$c = $ws.Range("D34")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '3.75'
$c.Style = $origStyle
$ws.Range("E34").Value = '  +3.54%  '
$c = $ws.Range("D35")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '7.15'
$c.Style = $origStyle
$ws.Range("E35").Value = '  +6.74%  '
$c = $ws.Range("D36")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '0.147'
$c.Style = $origStyle
$ws.Range("E36").Value = '  +13.20%  '
$ws.Range("B37").Value = 'Fetch.AI'
$ws.Range("C37").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$c = $ws.Range("D37")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '1.30'
$c.Style = $origStyle
$ws.Range("E37").Value = '  +4.84%  '
$ws.Range("B38").Value = 'EthereumClassic'
$ws.Range("C38").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$c = $ws.Range("D38")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '23.62'
$c.Style = $origStyle
$ws.Range("E38").Value = '  +6.85%  '
$ws.Range("E39").Value = '  +4.10%  '
$c = $ws.Range("D40")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '0.0891'
$c.Style = $origStyle
$ws.Range("E40").Value = '  +32.00%  '
$ws.Range("E41").Value = '  +0.18%  '
$ws.Range("E42").Value = '  -0.09%  '
$c = $ws.Range("D43")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '0.153'
$c.Style = $origStyle
$ws.Range("E43").Value = '  +16.81%  '
$c = $ws.Range("D44")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '0.403'
$c.Style = $origStyle
$ws.Range("E44").Value = '  +12.13%  '
$ws.Range("E45").Value = '  +0.01%  '
$ws.Range("E46").Value = '  +6.70%  '
$c = $ws.Range("D47")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '148.49'
$c.Style = $origStyle
$ws.Range("E47").Value = '  +0.51%  '
$c = $ws.Range("D48")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '43.89'
$c.Style = $origStyle
$ws.Range("E48").Value = '  +1.57%  '
$c = $ws.Range("D49")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '1.31'
$c.Style = $origStyle
$ws.Range("E49").Value = '  +9.18%  '
$ws.Range("B50").Value = 'Aave'
$ws.Range("C50").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$c = $ws.Range("D50")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '166.69'
$c.Style = $origStyle
$ws.Range("E50").Value = '  +6.04%  '
$ws.Range("B51").Value = 'Filecoin'
$ws.Range("C51").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$c = $ws.Range("D51")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '4.23'
$c.Style = $origStyle
$ws.Range("E51").Value = '  +8.01%  '
